$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (volume/issue number and week-covering date range) ---
$ws.Range("A8").Value = "Volume 30   Number  34"
$ws.Range("C9").Value = "Report Covering the Week  8/21/2023  Through  8/27/2023"

# --- Weekly crime-complaint statistics table (rows 14-30) ---
$ws.Range("C14").Value = 10
$ws.Range("D14").Value = 10
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 25
$ws.Range("G14").Value = 27
$ws.Range("H14").Value = -7.407407407407
$ws.Range("I14").Value = 263
$ws.Range("J14").Value = 293
$ws.Range("K14").Value = -10.238907849829
$ws.Range("L14").Value = -15.974440894568
$ws.Range("M14").Value = -25.495750708215
$ws.Range("N14").Value = -79.548989113530
$ws.Range("C15").Value = 25
$ws.Range("D15").Value = 33
$ws.Range("E15").Value = -24.242424242424
$ws.Range("F15").Value = 100
$ws.Range("G15").Value = 128
$ws.Range("H15").Value = -21.875
$ws.Range("I15").Value = 949
$ws.Range("J15").Value = 1084
$ws.Range("K15").Value = -12.453874538745
$ws.Range("L15").Value = -4.334677419354
$ws.Range("M15").Value = 9.710982658959
$ws.Range("N15").Value = -56.186518928901
$ws.Range("C16").Value = 379
$ws.Range("D16").Value = 358
$ws.Range("E16").Value = 5.865921787709
$ws.Range("F16").Value = 1431
$ws.Range("G16").Value = 1501
$ws.Range("H16").Value = -4.663557628247
$ws.Range("I16").Value = 10723
$ws.Range("J16").Value = 11381
$ws.Range("K16").Value = -5.781565767507
$ws.Range("L16").Value = 31.926673228346
$ws.Range("M16").Value = -11.358187980491
$ws.Range("N16").Value = -80.469546845402
$ws.Range("C17").Value = 505
$ws.Range("D17").Value = 540
$ws.Range("E17").Value = -6.481481481481
$ws.Range("F17").Value = 2203
$ws.Range("G17").Value = 2174
$ws.Range("H17").Value = 1.333946642134
$ws.Range("I17").Value = 18260
$ws.Range("J17").Value = 17333
$ws.Range("K17").Value = 5.348179772687
$ws.Range("L17").Value = 26.165964209217
$ws.Range("M17").Value = 61.293171981273
$ws.Range("N17").Value = -33.967381477597
$ws.Range("C18").Value = 235
$ws.Range("D18").Value = 356
$ws.Range("E18").Value = -33.988764044943
$ws.Range("F18").Value = 1058
$ws.Range("G18").Value = 1249
$ws.Range("H18").Value = -15.292233787029
$ws.Range("I18").Value = 9096
$ws.Range("J18").Value = 10209
$ws.Range("K18").Value = -10.90214516603
$ws.Range("L18").Value = 18.995290423861
$ws.Range("M18").Value = -23.363383604347
$ws.Range("N18").Value = -86.032370013206
$ws.Range("C19").Value = 1023
$ws.Range("D19").Value = 1014
$ws.Range("E19").Value = 0.887573964497
$ws.Range("F19").Value = 4195
$ws.Range("G19").Value = 4248
$ws.Range("H19").Value = -1.247645951035
$ws.Range("I19").Value = 32684
$ws.Range("J19").Value = 33424
$ws.Range("K19").Value = -2.213977979894
$ws.Range("L19").Value = 44.172915747684
$ws.Range("M19").Value = 36.342399466043
$ws.Range("N19").Value = -41.029157043879
$ws.Range("C20").Value = 353
$ws.Range("D20").Value = 297
$ws.Range("E20").Value = 18.855218855218
$ws.Range("F20").Value = 1364
$ws.Range("G20").Value = 1079
$ws.Range("H20").Value = 26.413345690454
$ws.Range("I20").Value = 10334
$ws.Range("J20").Value = 8669
$ws.Range("K20").Value = 19.206367516437
$ws.Range("L20").Value = 66.355441081777
$ws.Range("M20").Value = 53.711140859735
$ws.Range("N20").Value = -85.730067110386
$ws.Range("C21").Value = 2530
$ws.Range("D21").Value = 2608
$ws.Range("E21").Value = -2.990797546012
$ws.Range("F21").Value = 10376
$ws.Range("G21").Value = 10406
$ws.Range("H21").Value = -0.288295214299
$ws.Range("I21").Value = 82309
$ws.Range("J21").Value = 82393
$ws.Range("K21").Value = -0.101950408408
$ws.Range("L21").Value = 36.201019327508
$ws.Range("M21").Value = 22.483630952381
$ws.Range("N21").Value = -70.495711054474
$ws.Range("C22").Value = 38
$ws.Range("D22").Value = 48
$ws.Range("E22").Value = -20.833333333333
$ws.Range("F22").Value = 154
$ws.Range("G22").Value = 160
$ws.Range("H22").Value = -3.75
$ws.Range("I22").Value = 1423
$ws.Range("J22").Value = 1490
$ws.Range("K22").Value = -4.496644295302
$ws.Range("L22").Value = 43.158953722334
$ws.Range("M22").Value = 3.944485025566
$ws.Range("C23").Value = 128
$ws.Range("D23").Value = 113
$ws.Range("E23").Value = 13.274336283185
$ws.Range("F23").Value = 507
$ws.Range("G23").Value = 488
$ws.Range("H23").Value = 3.893442622950
$ws.Range("I23").Value = 4111
$ws.Range("J23").Value = 3959
$ws.Range("K23").Value = 3.839353372063
$ws.Range("L23").Value = 17.390062821245
$ws.Range("M23").Value = 51.418047882136
$ws.Range("C24").Value = 2267
$ws.Range("D24").Value = 2609
$ws.Range("E24").Value = -13.108470678420
$ws.Range("F24").Value = 8944
$ws.Range("G24").Value = 9724
$ws.Range("H24").Value = -8.021390374331
$ws.Range("I24").Value = 72303
$ws.Range("J24").Value = 74839
$ws.Range("K24").Value = -3.388607544194
$ws.Range("L24").Value = 38.175320580197
$ws.Range("M24").Value = 37.189533802629
$ws.Range("C25").Value = 814
$ws.Range("D25").Value = 788
$ws.Range("E25").Value = 3.299492385786
$ws.Range("F25").Value = 3446
$ws.Range("G25").Value = 3151
$ws.Range("H25").Value = 9.362107267534
$ws.Range("I25").Value = 28708
$ws.Range("J25").Value = 27331
$ws.Range("K25").Value = 5.038234971278
$ws.Range("L25").Value = 27.443842670691
$ws.Range("M25").Value = -5.819828095269
$ws.Range("C26").Value = 41
$ws.Range("D26").Value = 53
$ws.Range("E26").Value = -22.641509433962
$ws.Range("F26").Value = 177
$ws.Range("G26").Value = 197
$ws.Range("H26").Value = -10.152284263959
$ws.Range("I26").Value = 1592
$ws.Range("J26").Value = 1745
$ws.Range("K26").Value = -8.767908309455
$ws.Range("L26").Value = -1.179391682184
$ws.Range("C27").Value = 115
$ws.Range("D27").Value = 97
$ws.Range("E27").Value = 18.556701030927
$ws.Range("F27").Value = 422
$ws.Range("G27").Value = 429
$ws.Range("H27").Value = -1.631701631701
$ws.Range("I27").Value = 3470
$ws.Range("J27").Value = 3376
$ws.Range("K27").Value = 2.784360189573
$ws.Range("L27").Value = 11.791237113402
$ws.Range("C28").Value = 22
$ws.Range("D28").Value = 41
$ws.Range("E28").Value = -46.341463414634
$ws.Range("F28").Value = 93
$ws.Range("G28").Value = 131
$ws.Range("H28").Value = -29.007633587786
$ws.Range("I28").Value = 797
$ws.Range("J28").Value = 1107
$ws.Range("K28").Value = -28.003613369467
$ws.Range("L28").Value = -35.255889520714
$ws.Range("M28").Value = -32.799325463743
$ws.Range("N28").Value = -79.848293299620
$ws.Range("D29").Value = 31
$ws.Range("E29").Value = -29.032258064516
$ws.Range("F29").Value = 86
$ws.Range("H29").Value = -22.522522522522
$ws.Range("I29").Value = 677
$ws.Range("J29").Value = 916
$ws.Range("K29").Value = -26.091703056768
$ws.Range("L29").Value = -34.143968871595
$ws.Range("M29").Value = -31.129196337741
$ws.Range("N29").Value = -80.977802753582
$ws.Range("C30").Value = 1
$ws.Range("D30").Value = 15
$ws.Range("E30").Value = -93.333333333333
$ws.Range("F30").Value = 23
$ws.Range("G30").Value = 47
$ws.Range("H30").Value = -51.063829787234
$ws.Range("I30").Value = 309
$ws.Range("J30").Value = 451
$ws.Range("K30").Value = -31.485587583148
$ws.Range("L30").Value = -14.166666666666
